$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, "A").Value = "ECs"
$ws.Cells.Item(2, "B").Value = "Ghrl"
$ws.Cells.Item(2, "C").Value = "Ghsr"
$ws.Cells.Item(2, "D").Value = "ECs"
$ws.Cells.Item(2, "E").Value = 3
$ws.Cells.Item(2, "F").Value = 1
$ws.Cells.Item(2, "G").Value = 18.434123
$ws.Cells.Item(2, "H").Value = 55.302369
$ws.Cells.Item(2, "I").Value = 0.4617788569487251
$ws.Cells.Item(2, "J").Value = 0.4617788569487251
$ws.Cells.Item(2, "K").Value = 2
$ws.Cells.Item(2, "L").Value = 0.6666666666666666
$ws.Cells.Item(2, "M").Value = 0.4759383333333333
$ws.Cells.Item(2, "N").Value = 1.427815
$ws.Cells.Item(2, "O").Value = 0.1284588056245616
$ws.Cells.Item(2, "P").Value = 0.1284588056245616
$ws.Cells.Item(2, "Q").Value = 8.773505777081667
$ws.Cells.Item(2, "R").Value = 78.96155199373499
$ws.Cells.Item(2, "S").Value = 0.05931956042630853
$ws.Cells.Item(2, "T").Value = 0.05931956042630851

$ws.Cells.Item(3, "A").Value = "ECs"
$ws.Cells.Item(3, "B").Value = "Ghrl"
$ws.Cells.Item(3, "C").Value = "Ghsr"
$ws.Cells.Item(3, "D").Value = "FAPs"
$ws.Cells.Item(3, "E").Value = 3
$ws.Cells.Item(3, "F").Value = 1
$ws.Cells.Item(3, "G").Value = 18.434123
$ws.Cells.Item(3, "H").Value = 55.302369
$ws.Cells.Item(3, "I").Value = 0.4617788569487251
$ws.Cells.Item(3, "J").Value = 0.4617788569487251
$ws.Cells.Item(3, "K").Value = 3
$ws.Cells.Item(3, "L").Value = 1
$ws.Cells.Item(3, "M").Value = 2.371854333333333
$ws.Cells.Item(3, "N").Value = 7.115563
$ws.Cells.Item(3, "O").Value = 0.6401786816403544
$ws.Cells.Item(3, "P").Value = 0.6401786816403543
$ws.Cells.Item(3, "Q").Value = 43.72305451874966
$ws.Cells.Item(3, "R").Value = 393.507490668747
$ws.Cells.Item(3, "S").Value = 0.2956209798508246
$ws.Cells.Item(3, "T").Value = 0.2956209798508246

$ws.Cells.Item(4, "A").Value = "ECs"
$ws.Cells.Item(4, "B").Value = "Ghrl"
$ws.Cells.Item(4, "C").Value = "Ghsr"
$ws.Cells.Item(4, "D").Value = "M2"
$ws.Cells.Item(4, "E").Value = 3
$ws.Cells.Item(4, "F").Value = 1
$ws.Cells.Item(4, "G").Value = 18.434123
$ws.Cells.Item(4, "H").Value = 55.302369
$ws.Cells.Item(4, "I").Value = 0.4617788569487251
$ws.Cells.Item(4, "J").Value = 0.4617788569487251
$ws.Cells.Item(4, "K").Value = 1
$ws.Cells.Item(4, "L").Value = 0.3333333333333333
$ws.Cells.Item(4, "M").Value = 0.055189
$ws.Cells.Item(4, "N").Value = 0.165567
$ws.Cells.Item(4, "O").Value = 0.01489586470995318
$ws.Cells.Item(4, "P").Value = 0.01489586470995318
$ws.Cells.Item(4, "Q").Value = 1.017360814247
$ws.Cells.Item(4, "R").Value = 9.156247328223
$ws.Cells.Item(4, "S").Value = 0.00687859537902503
$ws.Cells.Item(4, "T").Value = 0.00687859537902503

$ws.Cells.Item(5, "A").Value = "ECs"
$ws.Cells.Item(5, "B").Value = "Ghrl"
$ws.Cells.Item(5, "C").Value = "Ghsr"
$ws.Cells.Item(5, "D").Value = "sCs"
$ws.Cells.Item(5, "E").Value = 3
$ws.Cells.Item(5, "F").Value = 1
$ws.Cells.Item(5, "G").Value = 18.434123
$ws.Cells.Item(5, "H").Value = 55.302369
$ws.Cells.Item(5, "I").Value = 0.4617788569487251
$ws.Cells.Item(5, "J").Value = 0.4617788569487251
$ws.Cells.Item(5, "K").Value = 3
$ws.Cells.Item(5, "L").Value = 1
$ws.Cells.Item(5, "M").Value = 0.8020063333333334
$ws.Cells.Item(5, "N").Value = 2.406019
$ws.Cells.Item(5, "O").Value = 0.2164666480251308
$ws.Cells.Item(5, "P").Value = 0.2164666480251308
$ws.Cells.Item(5, "Q").Value = 14.78428339544567
$ws.Cells.Item(5, "R").Value = 133.058550559011
$ws.Cells.Item(5, "S").Value = 0.09995972129256692
$ws.Cells.Item(5, "T").Value = 0.09995972129256689

$ws.Cells.Item(6, "A").Value = "FAPs"
$ws.Cells.Item(6, "B").Value = "Ghrl"
$ws.Cells.Item(6, "C").Value = "Ghsr"
$ws.Cells.Item(6, "D").Value = "ECs"
$ws.Cells.Item(6, "E").Value = 3
$ws.Cells.Item(6, "F").Value = 1
$ws.Cells.Item(6, "G").Value = 6.782643666666666
$ws.Cells.Item(6, "H").Value = 20.347931
$ws.Cells.Item(6, "I").Value = 0.1699067234977136
$ws.Cells.Item(6, "J").Value = 0.1699067234977136
$ws.Cells.Item(6, "K").Value = 2
$ws.Cells.Item(6, "L").Value = 0.6666666666666666
$ws.Cells.Item(6, "M").Value = 0.4759383333333333
$ws.Cells.Item(6, "N").Value = 1.427815
$ws.Cells.Item(6, "O").Value = 0.1284588056245616
$ws.Cells.Item(6, "P").Value = 0.1284588056245616
$ws.Cells.Item(6, "Q").Value = 3.228120122307222
$ws.Cells.Item(6, "R").Value = 29.053081100765
$ws.Cells.Item(6, "S").Value = 0.02182601476809893
$ws.Cells.Item(6, "T").Value = 0.02182601476809893

$ws.Cells.Item(7, "A").Value = "FAPs"
$ws.Cells.Item(7, "B").Value = "Ghrl"
$ws.Cells.Item(7, "C").Value = "Ghsr"
$ws.Cells.Item(7, "D").Value = "FAPs"
$ws.Cells.Item(7, "E").Value = 3
$ws.Cells.Item(7, "F").Value = 1
$ws.Cells.Item(7, "G").Value = 6.782643666666666
$ws.Cells.Item(7, "H").Value = 20.347931
$ws.Cells.Item(7, "I").Value = 0.1699067234977136
$ws.Cells.Item(7, "J").Value = 0.1699067234977136
$ws.Cells.Item(7, "K").Value = 3
$ws.Cells.Item(7, "L").Value = 1
$ws.Cells.Item(7, "M").Value = 2.371854333333333
$ws.Cells.Item(7, "N").Value = 7.115563
$ws.Cells.Item(7, "O").Value = 0.6401786816403544
$ws.Cells.Item(7, "P").Value = 0.6401786816403543
$ws.Cells.Item(7, "Q").Value = 16.08744277223922
$ws.Cells.Item(7, "R").Value = 144.786984950153
$ws.Cells.Item(7, "S").Value = 0.1087706622505985
$ws.Cells.Item(7, "T").Value = 0.1087706622505985

$ws.Cells.Item(8, "A").Value = "FAPs"
$ws.Cells.Item(8, "B").Value = "Ghrl"
$ws.Cells.Item(8, "C").Value = "Ghsr"
$ws.Cells.Item(8, "D").Value = "M2"
$ws.Cells.Item(8, "E").Value = 3
$ws.Cells.Item(8, "F").Value = 1
$ws.Cells.Item(8, "G").Value = 6.782643666666666
$ws.Cells.Item(8, "H").Value = 20.347931
$ws.Cells.Item(8, "I").Value = 0.1699067234977136
$ws.Cells.Item(8, "J").Value = 0.1699067234977136
$ws.Cells.Item(8, "K").Value = 1
$ws.Cells.Item(8, "L").Value = 0.3333333333333333
$ws.Cells.Item(8, "M").Value = 0.055189
$ws.Cells.Item(8, "N").Value = 0.165567
$ws.Cells.Item(8, "O").Value = 0.01489586470995318
$ws.Cells.Item(8, "P").Value = 0.01489586470995318
$ws.Cells.Item(8, "Q").Value = 0.3743273213196666
$ws.Cells.Item(8, "R").Value = 3.368945891877
$ws.Cells.Item(8, "S").Value = 0.002530907566533364
$ws.Cells.Item(8, "T").Value = 0.002530907566533364

$ws.Cells.Item(9, "A").Value = "FAPs"
$ws.Cells.Item(9, "B").Value = "Ghrl"
$ws.Cells.Item(9, "C").Value = "Ghsr"
$ws.Cells.Item(9, "D").Value = "sCs"
$ws.Cells.Item(9, "E").Value = 3
$ws.Cells.Item(9, "F").Value = 1
$ws.Cells.Item(9, "G").Value = 6.782643666666666
$ws.Cells.Item(9, "H").Value = 20.347931
$ws.Cells.Item(9, "I").Value = 0.1699067234977136
$ws.Cells.Item(9, "J").Value = 0.1699067234977136
$ws.Cells.Item(9, "K").Value = 3
$ws.Cells.Item(9, "L").Value = 1
$ws.Cells.Item(9, "M").Value = 0.8020063333333334
$ws.Cells.Item(9, "N").Value = 2.406019
$ws.Cells.Item(9, "O").Value = 0.2164666480251308
$ws.Cells.Item(9, "P").Value = 0.2164666480251308
$ws.Cells.Item(9, "Q").Value = 5.439723177409888
$ws.Cells.Item(9, "R").Value = 48.957508596689
$ws.Cells.Item(9, "S").Value = 0.03677913891248279
$ws.Cells.Item(9, "T").Value = 0.03677913891248279

$ws.Cells.Item(10, "A").Value = "M2"
$ws.Cells.Item(10, "B").Value = "Ghrl"
$ws.Cells.Item(10, "C").Value = "Ghsr"
$ws.Cells.Item(10, "D").Value = "ECs"
$ws.Cells.Item(10, "E").Value = 3
$ws.Cells.Item(10, "F").Value = 1
$ws.Cells.Item(10, "G").Value = 11.24197566666667
$ws.Cells.Item(10, "H").Value = 33.725927
$ws.Cells.Item(10, "I").Value = 0.2816139760594369
$ws.Cells.Item(10, "J").Value = 0.2816139760594369
$ws.Cells.Item(10, "K").Value = 2
$ws.Cells.Item(10, "L").Value = 0.6666666666666666
$ws.Cells.Item(10, "M").Value = 0.4759383333333333
$ws.Cells.Item(10, "N").Value = 1.427815
$ws.Cells.Item(10, "O").Value = 0.1284588056245616
$ws.Cells.Item(10, "P").Value = 0.1284588056245616
$ws.Cells.Item(10, "Q").Value = 5.350487162167222
$ws.Cells.Item(10, "R").Value = 48.15438445950499
$ws.Cells.Item(10, "S").Value = 0.03617579501177916
$ws.Cells.Item(10, "T").Value = 0.03617579501177915

$ws.Cells.Item(11, "A").Value = "M2"
$ws.Cells.Item(11, "B").Value = "Ghrl"
$ws.Cells.Item(11, "C").Value = "Ghsr"
$ws.Cells.Item(11, "D").Value = "FAPs"
$ws.Cells.Item(11, "E").Value = 3
$ws.Cells.Item(11, "F").Value = 1
$ws.Cells.Item(11, "G").Value = 11.24197566666667
$ws.Cells.Item(11, "H").Value = 33.725927
$ws.Cells.Item(11, "I").Value = 0.2816139760594369
$ws.Cells.Item(11, "J").Value = 0.2816139760594369
$ws.Cells.Item(11, "K").Value = 3
$ws.Cells.Item(11, "L").Value = 1
$ws.Cells.Item(11, "M").Value = 2.371854333333333
$ws.Cells.Item(11, "N").Value = 7.115563
$ws.Cells.Item(11, "O").Value = 0.6401786816403544
$ws.Cells.Item(11, "P").Value = 0.6401786816403543
$ws.Cells.Item(11, "Q").Value = 26.66432870021122
$ws.Cells.Item(11, "R").Value = 239.978958301901
$ws.Cells.Item(11, "S").Value = 0.1802832639252286
$ws.Cells.Item(11, "T").Value = 0.1802832639252286

$ws.Cells.Item(12, "A").Value = "M2"
$ws.Cells.Item(12, "B").Value = "Ghrl"
$ws.Cells.Item(12, "C").Value = "Ghsr"
$ws.Cells.Item(12, "D").Value = "M2"
$ws.Cells.Item(12, "E").Value = 3
$ws.Cells.Item(12, "F").Value = 1
$ws.Cells.Item(12, "G").Value = 11.24197566666667
$ws.Cells.Item(12, "H").Value = 33.725927
$ws.Cells.Item(12, "I").Value = 0.2816139760594369
$ws.Cells.Item(12, "J").Value = 0.2816139760594369
$ws.Cells.Item(12, "K").Value = 1
$ws.Cells.Item(12, "L").Value = 0.3333333333333333
$ws.Cells.Item(12, "M").Value = 0.055189
$ws.Cells.Item(12, "N").Value = 0.165567
$ws.Cells.Item(12, "O").Value = 0.01489586470995318
$ws.Cells.Item(12, "P").Value = 0.01489586470995318
$ws.Cells.Item(12, "Q").Value = 0.6204333950676666
$ws.Cells.Item(12, "R").Value = 5.583900555609
$ws.Cells.Item(12, "S").Value = 0.004194883687813364
$ws.Cells.Item(12, "T").Value = 0.004194883687813364

$ws.Cells.Item(13, "A").Value = "M2"
$ws.Cells.Item(13, "B").Value = "Ghrl"
$ws.Cells.Item(13, "C").Value = "Ghsr"
$ws.Cells.Item(13, "D").Value = "sCs"
$ws.Cells.Item(13, "E").Value = 3
$ws.Cells.Item(13, "F").Value = 1
$ws.Cells.Item(13, "G").Value = 11.24197566666667
$ws.Cells.Item(13, "H").Value = 33.725927
$ws.Cells.Item(13, "I").Value = 0.2816139760594369
$ws.Cells.Item(13, "J").Value = 0.2816139760594369
$ws.Cells.Item(13, "K").Value = 3
$ws.Cells.Item(13, "L").Value = 1
$ws.Cells.Item(13, "M").Value = 0.8020063333333334
$ws.Cells.Item(13, "N").Value = 2.406019
$ws.Cells.Item(13, "O").Value = 0.2164666480251308
$ws.Cells.Item(13, "P").Value = 0.2164666480251308
$ws.Cells.Item(13, "Q").Value = 9.01613568384589
$ws.Cells.Item(13, "R").Value = 81.145221154613
$ws.Cells.Item(13, "S").Value = 0.06096003343461574
$ws.Cells.Item(13, "T").Value = 0.06096003343461574

$ws.Cells.Item(14, "A").Value = "sCs"
$ws.Cells.Item(14, "B").Value = "Ghrl"
$ws.Cells.Item(14, "C").Value = "Ghsr"
$ws.Cells.Item(14, "D").Value = "ECs"
$ws.Cells.Item(14, "E").Value = 3
$ws.Cells.Item(14, "F").Value = 1
$ws.Cells.Item(14, "G").Value = 3.461065
$ws.Cells.Item(14, "H").Value = 10.383195
$ws.Cells.Item(14, "I").Value = 0.08670044349412441
$ws.Cells.Item(14, "J").Value = 0.0867004434941244
$ws.Cells.Item(14, "K").Value = 2
$ws.Cells.Item(14, "L").Value = 0.6666666666666666
$ws.Cells.Item(14, "M").Value = 0.4759383333333333
$ws.Cells.Item(14, "N").Value = 1.427815
$ws.Cells.Item(14, "O").Value = 0.1284588056245616
$ws.Cells.Item(14, "P").Value = 0.1284588056245616
$ws.Cells.Item(14, "Q").Value = 1.647253507658333
$ws.Cells.Item(14, "R").Value = 14.825281568925
$ws.Cells.Item(14, "S").Value = 0.01113743541837502
$ws.Cells.Item(14, "T").Value = 0.01113743541837501

$ws.Cells.Item(15, "A").Value = "sCs"
$ws.Cells.Item(15, "B").Value = "Ghrl"
$ws.Cells.Item(15, "C").Value = "Ghsr"
$ws.Cells.Item(15, "D").Value = "FAPs"
$ws.Cells.Item(15, "E").Value = 3
$ws.Cells.Item(15, "F").Value = 1
$ws.Cells.Item(15, "G").Value = 3.461065
$ws.Cells.Item(15, "H").Value = 10.383195
$ws.Cells.Item(15, "I").Value = 0.08670044349412441
$ws.Cells.Item(15, "J").Value = 0.0867004434941244
$ws.Cells.Item(15, "K").Value = 3
$ws.Cells.Item(15, "L").Value = 1
$ws.Cells.Item(15, "M").Value = 2.371854333333333
$ws.Cells.Item(15, "N").Value = 7.115563
$ws.Cells.Item(15, "O").Value = 0.6401786816403544
$ws.Cells.Item(15, "P").Value = 0.6401786816403543
$ws.Cells.Item(15, "Q").Value = 8.209142018198332
$ws.Cells.Item(15, "R").Value = 73.88227816378499
$ws.Cells.Item(15, "S").Value = 0.05550377561370261
$ws.Cells.Item(15, "T").Value = 0.05550377561370259

$ws.Cells.Item(16, "A").Value = "sCs"
$ws.Cells.Item(16, "B").Value = "Ghrl"
$ws.Cells.Item(16, "C").Value = "Ghsr"
$ws.Cells.Item(16, "D").Value = "M2"
$ws.Cells.Item(16, "E").Value = 3
$ws.Cells.Item(16, "F").Value = 1
$ws.Cells.Item(16, "G").Value = 3.461065
$ws.Cells.Item(16, "H").Value = 10.383195
$ws.Cells.Item(16, "I").Value = 0.08670044349412441
$ws.Cells.Item(16, "J").Value = 0.0867004434941244
$ws.Cells.Item(16, "K").Value = 1
$ws.Cells.Item(16, "L").Value = 0.3333333333333333
$ws.Cells.Item(16, "M").Value = 0.055189
$ws.Cells.Item(16, "N").Value = 0.165567
$ws.Cells.Item(16, "O").Value = 0.01489586470995318
$ws.Cells.Item(16, "P").Value = 0.01489586470995318
$ws.Cells.Item(16, "Q").Value = 0.191012716285
$ws.Cells.Item(16, "R").Value = 1.719114446565
$ws.Cells.Item(16, "S").Value = 0.001291478076581417
$ws.Cells.Item(16, "T").Value = 0.001291478076581417

$ws.Cells.Item(17, "A").Value = "sCs"
$ws.Cells.Item(17, "B").Value = "Ghrl"
$ws.Cells.Item(17, "C").Value = "Ghsr"
$ws.Cells.Item(17, "D").Value = "sCs"
$ws.Cells.Item(17, "E").Value = 3
$ws.Cells.Item(17, "F").Value = 1
$ws.Cells.Item(17, "G").Value = 3.461065
$ws.Cells.Item(17, "H").Value = 10.383195
$ws.Cells.Item(17, "I").Value = 0.08670044349412441
$ws.Cells.Item(17, "J").Value = 0.0867004434941244
$ws.Cells.Item(17, "K").Value = 3
$ws.Cells.Item(17, "L").Value = 1
$ws.Cells.Item(17, "M").Value = 0.8020063333333334
$ws.Cells.Item(17, "N").Value = 2.406019
$ws.Cells.Item(17, "O").Value = 0.2164666480251308
$ws.Cells.Item(17, "P").Value = 0.2164666480251308
$ws.Cells.Item(17, "Q").Value = 2.775796050078333
$ws.Cells.Item(17, "R").Value = 24.982164450705
$ws.Cells.Item(17, "S").Value = 0.01876775438546537
$ws.Cells.Item(17, "T").Value = 0.01876775438546537
